$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.934329628944397
$ws.Range("B1").Value = 1.897062182426453
$ws.Range("C1").Value = 4.374848365783691
$ws.Range("D1").Value = 3.366589546203613
$ws.Range("E1").Value = 1.457905292510986
